# Ran code for averaged intensities on spiral schemes.
# "Gaussian-Quadrature" now leads a second scheme block, followed by three new
# "Spiral-90deg-*" sampling runs; the previously-existing NoRotation / Rotation /
# HexGrid rows are kept (re-seated further down the table) with their original data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 11-19 sit beyond today's used range (A1:M16), so first stamp them with the
# same bordered / bold / centered style already used by column A's HKL-index cells
# (copy format only -- this mirrors rows 11, 17, 18, 19 which are brand new, and
# is a no-op in substance for rows that already carry it).
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A10:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.18355117195812
$ws.Range("D10").Value = 0.5960752793557353
$ws.Range("E10").Value = 1.047915988402299
$ws.Range("F10").Value = 1.18355117195812
$ws.Range("G10").Value = 0.7931908076294506
$ws.Range("H10").Value = 1.124909215325272
$ws.Range("I10").Value = 1.092202312283082
$ws.Range("J10").Value = 0.5960752793557353
$ws.Range("K10").Value = 0.8219956338790171
$ws.Range("L10").Value = 1.002773402918569
$ws.Range("M10").Value = 0.9729741291589932

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.109352432680098
$ws.Range("D11").Value = 0.7372935530304323
$ws.Range("E11").Value = 1.013445883063445
$ws.Range("F11").Value = 1.109352432680098
$ws.Range("G11").Value = 0.9140917390995468
$ws.Range("H11").Value = 0.9705272221040968
$ws.Range("I11").Value = 1.042395897866687
$ws.Range("J11").Value = 0.7372935530304323
$ws.Range("K11").Value = 0.8753697180469387
$ws.Range("L11").Value = 0.9923610753635181
$ws.Range("M11").Value = 0.9645177879740511

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.107420528287933
$ws.Range("D12").Value = 0.7383581386454428
$ws.Range("E12").Value = 1.013908618270609
$ws.Range("F12").Value = 1.107420528287933
$ws.Range("G12").Value = 0.9147400249322111
$ws.Range("H12").Value = 0.9709702401372768
$ws.Range("I12").Value = 1.042070417114944
$ws.Range("J12").Value = 0.7383581386454428
$ws.Range("K12").Value = 0.8761333784580256
$ws.Range("L12").Value = 0.9917769533729792
$ws.Range("M12").Value = 0.964577994564736

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.109339627063266
$ws.Range("D13").Value = 0.7371847532975258
$ws.Range("E13").Value = 1.013537641279112
$ws.Range("F13").Value = 1.109339627063266
$ws.Range("G13").Value = 0.9141252101395791
$ws.Range("H13").Value = 0.9702996556195644
$ws.Range("I13").Value = 1.04226144711371
$ws.Range("J13").Value = 0.7371847532975258
$ws.Range("K13").Value = 0.8753611972883191
$ws.Range("L13").Value = 0.9923504121757925
$ws.Range("M13").Value = 0.9644580557521264

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.8039439999999991
$ws.Range("D14").Value = 0.4209960000000003
$ws.Range("E14").Value = 1.293879999999999
$ws.Range("F14").Value = 0.8039439999999991
$ws.Range("G14").Value = 0.5486320000000005
$ws.Range("H14").Value = 2.065675999999999
$ws.Range("I14").Value = 1.201939999999998
$ws.Range("J14").Value = 0.4209960000000003
$ws.Range("K14").Value = 0.8574379999999995
$ws.Range("L14").Value = 0.8306909999999993
$ws.Range("M14").Value = 1.055844666666666

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.65
$ws.Range("D15").Value = 0.1
$ws.Range("E15").Value = 1.5
$ws.Range("F15").Value = 0.65
$ws.Range("G15").Value = 0.2195875
$ws.Range("H15").Value = 2.91
$ws.Range("I15").Value = 1.35
$ws.Range("J15").Value = 0.1
$ws.Range("K15").Value = 0.8
$ws.Range("L15").Value = 0.725
$ws.Range("M15").Value = 1.121597916666667

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.8045109870592011
$ws.Range("D16").Value = 0.4685303427072001
$ws.Range("E16").Value = 1.2889961391104
$ws.Range("F16").Value = 0.8045109870592011
$ws.Range("G16").Value = 0.5439500617728007
$ws.Range("H16").Value = 2.089979913420804
$ws.Range("I16").Value = 1.197365059788797
$ws.Range("J16").Value = 0.4685303427072001
$ws.Range("K16").Value = 0.8787632409088
$ws.Range("L16").Value = 0.8416371139840005
$ws.Range("M16").Value = 1.065555417309867

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9949238975122726
$ws.Range("D17").Value = 0.9949169638402917
$ws.Range("E17").Value = 0.9943649668699153
$ws.Range("F17").Value = 0.9949238975122726
$ws.Range("G17").Value = 0.9916734764235375
$ws.Range("H17").Value = 0.9951561430000584
$ws.Range("I17").Value = 0.99493714788777
$ws.Range("J17").Value = 0.9949169638402917
$ws.Range("K17").Value = 0.9946409653551035
$ws.Range("L17").Value = 0.9947824314336881
$ws.Range("M17").Value = 0.9943287659223076

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9191403110875501
$ws.Range("D18").Value = 1.060380223990361
$ws.Range("E18").Value = 1.00861947176323
$ws.Range("F18").Value = 0.9191403110875501
$ws.Range("G18").Value = 1.022333406840271
$ws.Range("H18").Value = 1.013690995402253
$ws.Range("I18").Value = 0.9826149833148545
$ws.Range("J18").Value = 1.060380223990361
$ws.Range("K18").Value = 1.034499847876796
$ws.Range("L18").Value = 0.976820079482173
$ws.Range("M18").Value = 1.001129898733087

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9724300142247012
$ws.Range("D19").Value = 1.163659400803345
$ws.Range("E19").Value = 0.949661544349785
$ws.Range("F19").Value = 0.9724300142247012
$ws.Range("G19").Value = 1.08407045929631
$ws.Range("H19").Value = 0.8719854711094791
$ws.Range("I19").Value = 0.9506252903884476
$ws.Range("J19").Value = 1.163659400803345
$ws.Range("K19").Value = 1.056660472576565
$ws.Range("L19").Value = 1.014545243400633
$ws.Range("M19").Value = 0.9987386966953444

